# Auto-generated edit script: updates cryptos.xlsx price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.414.55"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3
$ws.Range("D3").Value = "1.911.74"
$ws.Range("E3").Value = "  +1.08%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.98%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.85"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$ws.Range("E6").Value = "  +0.92%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4812"
$ws.Range("E7").Value = "  +1.40%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4070"
$ws.Range("E8").Value = "  +0.72%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08204"
$ws.Range("E9").Value = "  +2.34%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.019"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.44"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12
$ws.Range("D12").Value = "1.918.05"
$ws.Range("E12").Value = "  +1.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.032"
$ws.Range("E13").Value = "  +1.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.210"
$ws.Range("E14").Value = "  +2.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.94"
$ws.Range("E15").Value = "  +1.92%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06806"
$ws.Range("E16").Value = "  +2.50%  "

# Row 17
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +0.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  +0.75%  "

# Row 21
$ws.Range("D21").Value = "29.444.22"
$ws.Range("E21").Value = "  +0.55%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.616"
$ws.Range("E22").Value = "  +2.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.72"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("E24").Value = "  +1.19%  "

# Row 25
$ws.Range("D25").Value = "2.151.00"
$ws.Range("E25").Value = "  +1.70%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.590"
$ws.Range("E26").Value = "  +10.73%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.96"
$ws.Range("E27").Value = "  +1.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  +2.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.104"
$ws.Range("E29").Value = "  +1.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.01"
$ws.Range("E30").Value = "  +2.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.019"
$ws.Range("E31").Value = "  -0.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09549"
$ws.Range("E32").Value = "  +1.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.549"
$ws.Range("E33").Value = "  +3.75%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.559"
$ws.Range("E34").Value = "  +0.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.366"
$ws.Range("E35").Value = "  -0.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02282"
$ws.Range("E36").Value = "  +1.60%  "

# Row 37
$ws.Range("E37").Value = "  +1.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.178"
$ws.Range("E38").Value = "  +0.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5970"
$ws.Range("E39").Value = "  +2.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.037"
$ws.Range("E40").Value = "  +1.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.78"
$ws.Range("E41").Value = "  +7.37%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1850"
$ws.Range("E42").Value = "  +0.92%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.418"
$ws.Range("E43").Value = "  +1.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.245"
$ws.Range("E44").Value = "  -3.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07585"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.38"
$ws.Range("E46").Value = "  +1.41%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5568"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.949"
$ws.Range("E48").Value = "  +2.17%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.22"
$ws.Range("E49").Value = "  +3.59%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.422"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.42"
$ws.Range("E51").Value = "  +1.68%  "

